# Reformatted recipes and ingredients + started index rework
#
# The only functional content edit in this commit is to cell C16 on
# Sheet1: it held the number 40, and is replaced with the text "Tomato 3".
# (The SUM formula in J11 references C16, so once C16 becomes text it is
# excluded from the sum and the cached total drops from 124 to 84 — Excel
# recalculates this automatically.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C16").Value = "Tomato 3"
